$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 4 ("login.spec.ts" / "test login page22" / ...) - rows below shift up
$ws.Rows(4).Delete()

# Fix the wording on row 3 ("test login page" -> "login page")
$ws.Range("B3").Value = "login page"

# Add the two new header columns
$ws.Range("I1").Value = "User Story"
$ws.Range("J1").Value = "Test Case ID"

# Populate the new "User Story" / "Test Case ID" columns for each data row.
# Test Case ID values look numeric ("666", "777", "1234") but must stay text,
# so format those cells as Text before writing, then drop back to the
# Normal style (keeps the text type without leaving a custom number format
# attached to the cell).
$ws.Range("I2").Value = "https://jira.app.syfbank.com/browse/SAK-1over"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "666"
$ws.Range("J2").Style = "Normal"

$ws.Range("I3").Value = "https://jira.app.syfbank.com/browse/SAK-143"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "777"
$ws.Range("J3").Style = "Normal"

$ws.Range("I4").Value = "https://jira.app.syfbank.com/browse/SAK-123"
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "1234"
$ws.Range("J4").Style = "Normal"

# Row 5 (reqres.spec.ts) has no User Story / Test Case ID - left blank.

# Widen the two new columns to match the rest of the sheet's custom widths
$ws.Columns("I:J").ColumnWidth = 19.166666666666668
